# ---------------------------------------------------------------------------
# ReadMe_Configure Rest API Projects.docx
#   1) Mark every inline picture's run as "no proofing" (<w:noProof/> in the
#      run's rPr) -- matches what Word stamps on a run holding a drawing.
#   2) Append the new Q&A content (JSON request/response instructions) after
#      the existing closing paragraph, moving the trailing "_GoBack" bookmark
#      to sit in front of the brand-new final run.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Flag every InlineShape's run as noProof -----------------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = 1
}

# --- 2) Append the new Q&A paragraphs ---------------------------------------
# The "_GoBack" bookmark currently sits right at the very end of the last
# paragraph. Drop it first so the new content lands after the existing text
# without Word pinning the bookmark to the old location; it gets re-created
# in its correct spot (just before the final new run) as part of the XML
# below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$newParagraphsXml = (
    '<w:p xmlns:w="' + $w + '"/>' +
    '<w:p xmlns:w="' + $w + '">' +
        '<w:r><w:t>Q1. How to make this Project work for Json Request and response?</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $w + '"/>' +
    '<w:p xmlns:w="' + $w + '">' +
        '<w:r><w:t xml:space="preserve">Add </w:t></w:r>' +
        '<w:r><w:t>genson-1.3.jar</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
        '<w:r><w:t>jersey-json-1.8.jar</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> jar files </w:t></w:r>' +
        '<w:r><w:t>in WEB-INF - &gt; lib</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $w + '"/>' +
    '<w:p xmlns:w="' + $w + '">' +
        '<w:r><w:t xml:space="preserve">For Get </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">/Post </w:t></w:r>' +
        '<w:r><w:t xml:space="preserve">Request </w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $w + '"/>' +
    '<w:p xmlns:w="' + $w + '">' +
        '<w:r><w:t>@Produces(MediaType.APPLICATION_JSON)</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $w + '">' +
        '<w:r><w:t>Patient getPatient(@PathParam("id") String id);</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $w + '"/>' +
    '<w:p xmlns:w="' + $w + '">' +
        '<w:r><w:t>@Consumes(MediaType.APPLICATION_JSON)</w:t></w:r>' +
    '</w:p>' +
    '<w:p xmlns:w="' + $w + '">' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:t>Response addPatient(Patient patient);</w:t></w:r>' +
    '</w:p>'
)

# Build a brand-new collapsed Range at the very end of the document (using
# Range(pos, pos) rather than Range.Collapse, which this host does not fully
# honour for subsequent InsertXML calls) and splice the OOXML in after the
# existing last paragraph.
$lastEnd = $d.Paragraphs.Last.Range.End
$insertionPoint = $d.Range($lastEnd, $lastEnd)
$insertionPoint.InsertXML($newParagraphsXml)
